# Gallery_PanelNode_Repeaters_MZX_Panels.xlsx
# "Test data added for Italy"
#
# - Adds a new "Italy" worksheet (as the last tab) with the same layout as
#   the other full country sheets (Germany/Belgium/Czech/Portugal): header
#   block in rows 1-7, then one row per repeater model in rows 8-21.
# - Market name -> "Italy Market", User Story code -> "NGC-3145/T2454/T2453".
# - Leaves a selection "bookmark" on the Germany sheet (A8:A21, where the
#   repeater list starts) and clears/resets Slovakia's selection to a
#   whole-sheet selection now that it is no longer the active tab.
# - Italy becomes the active tab, with the selection sitting just below the
#   pasted data (row 22).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# --- Germany: leave the selection on the repeater list (A8:A21) ---------
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
[void]$germany.Range("A8:A21").Select()

# --- Slovakia: select the whole sheet (it will no longer be the active tab) ---
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Activate()
[void]$slovakia.Cells.Select()

# --- Build "Italy" from a copy of "Slovakia" (same column widths/styles) ---
$slovakia.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Slovakia's sheet only lists 9 repeater models; the other full country
# sheets list 11 (it also has P32AR / P32DR). Insert the two missing rows
# above "PR1DS" and copy that row's formatting onto them.
$italy.Rows.Item(16).Resize(2).Insert()
$italy.Range("A18").Copy()
$italy.Range("A16:A17").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$italy.Range("A16").Value = "P32AR"
$italy.Range("A17").Value = "P32DR"

# Shared strings are appended in the order they're first written, so set the
# User Story code before the market name to match the authored string order.
$italy.Range("B4").Value = "NGC-3145/T2454/T2453"
$italy.Range("B2").Value = "Italy Market"

# Italy is the newly active tab, selection resting one row below the data.
$italy.Activate()
[void]$italy.Rows.Item(22).Select()
